$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-14 (Uke, Dag, Navn, Lunsj)
$data = @(
    @(3,  "Mandag",  "Susanne", "Sushi"),
    @(5,  "Onsdag",  "Susanne", "Sushi"),
    @(7,  "Fredag",  "Susanne", "Poteter og fisk"),
    @(11, "Fredag",  "Susanne", "Lasagne"),
    @(15, "Torsdag", "Susanne", "Brødskiver med pålegg"),
    @(21, "Tirsdag", "Susanne", "Sushi"),
    @(23, "Fredag",  "Susanne", "Salat"),
    @(25, "Fredag",  "Susanne", "Pasta bolognese"),
    @(35, "Mandag",  "Susanne", "Ris og kylling"),
    @(37, "Mandag",  "Susanne", "Pasta bolognese"),
    @(39, "Torsdag", "Susanne", "Salat"),
    @(45, "Tirsdag", "Susanne", "Pølser"),
    @(49, "Tirsdag", "Susanne", "Lasagne")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
